# "Generate Report for Archive"
# - Update localization status text from "Ready for handoff" to "In Translation"
#   everywhere it is used (Overview!E2:F2, zh-cn!C2, de-de!C2).
# - Narrow the now-shorter "Status" columns (Overview E:F, zh-cn C, de-de C)
#   to match the re-generated report's best-fit width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Replace the status text in every cell that currently shows it, so the
# workbook ends up with a single shared string shared by all occurrences.
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Shrink the Status columns to fit the new, shorter text.
$wsOverview.Range("E1").ColumnWidth = 12.43
$wsOverview.Range("F1").ColumnWidth = 12.43
$wsZhCn.Range("C1").ColumnWidth = 12.43
$wsDeDe.Range("C1").ColumnWidth = 12.43
